$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.852.02"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "2.244.51"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.90%  "
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.481"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0791"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "2.588.89"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "2.231.53"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").Value = "40.783.85"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "238.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +3.76%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.09%  "
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0727"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.105"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.87%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.115"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("E41").Value = "  +5.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").Value = "2.094.80"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.00%  "
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.97%  "
$ws.Range("E48").Value = "  -14.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D50").Value = "2.460.20"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("E51").Value = "  +3.55%  "
